$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2€")

# Copy the row-21 look & feel down to row 22 first (borders, fonts, number formats, alignment)
$ws.Range("A21:I21").Copy() | Out-Null
$ws.Range("A22:I22").PasteSpecial(-4122) | Out-Null

# New row of data (2025, County of Carladès variety)
$ws.Range("A22").Value = 2025
$ws.Range("B22").Value = "County of Carladès"
$ws.Range("C22").ClearContents() | Out-Null
$ws.Range("D22").Value = "Obv: With mint symbol - Cornucopia"
$ws.Range("E22").Value = "Rev: new map of Europe"
$ws.Range("F22").Value = "Obv: Mint main engraving Symbol - Square"

# "15.000" looks numeric under the en-US locale, force it to stay text like the other mintage cells
$ws.Range("G22").Formula = '="15.000"'
$ws.Range("G22").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4163) | Out-Null

$ws.Range("H22").Value = 1
$ws.Range("I22").Formula = '=IF(OR(AND(H22>1,H22<>"-")),"Can exchange","")'

# Highlight the new mintage value with a custom purple font/fill to flag the variety
$ws.Range("H22").Font.Color = 10498160
$ws.Range("H22").Interior.Color = 10498160
$ws.Range("H22").HorizontalAlignment = -4108

$ws.Range("G10").Select()
$ws.Range("L28").Select()

$wb.Worksheets.Item("Links").Range("B7").Select()
